$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G4").Value = 101
$ws.Range("H4").Value = 1139
$ws.Range("I4").Value = 1010
$ws.Range("J4").Value = 1102
$ws.Range("Q4").Value = 722
